# Update Tgfb1-Sdc2 LR-pair TPM data (YoungD4) with recomputed values.
# Source data now includes "Resolving-Mac" as a Target cluster in addition
# to being a Sending cluster, so the sheet grows from 15 to 20 data rows
# (5 senders x 4 targets) and every numeric column is refreshed with the
# new TPM-based computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster / Ligand symbol / Receptor symbol / Target cluster, then
# the 16 numeric metric columns (E..T), for each of the 20 rows (rows 2..21).
$rows = @(
    @("ECs", "Tgfb1", "Sdc2", "ECs", 3, 1, 47.86240033333333, 143.587201, 0.1228118231805696, 0.1228118231805696, 3, 1, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 87.42383312725534, 786.814498145298, 0.003433048429305734, 0.003433048429305734),
    @("ECs", "Tgfb1", "Sdc2", "FAPs", 3, 1, 47.86240033333333, 143.587201, 0.1228118231805696, 0.1228118231805696, 3, 1, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 2119.82856655322, 19078.45709897898, 0.08324359468670023, 0.08324359468670021),
    @("ECs", "Tgfb1", "Sdc2", "MuSCs", 3, 1, 47.86240033333333, 143.587201, 0.1228118231805696, 0.1228118231805696, 3, 1, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 914.5963539493565, 8231.367185544208, 0.03591530437477131, 0.0359153043747713),
    @("ECs", "Tgfb1", "Sdc2", "Resolving-Mac", 3, 1, 47.86240033333333, 143.587201, 0.1228118231805696, 0.1228118231805696, 2, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 5.599214811261889, 50.39293330135699, 0.0002198756897923667, 0.0002198756897923667),
    @("FAPs", "Tgfb1", "Sdc2", "ECs", 3, 1, 18.65566266666667, 55.966988, 0.04786922362394307, 0.04786922362394307, 3, 1, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 34.07579913440266, 306.682192209624, 0.001338123307009605, 0.001338123307009605),
    @("FAPs", "Tgfb1", "Sdc2", "FAPs", 3, 1, 18.65566266666667, 55.966988, 0.04786922362394307, 0.04786922362394307, 3, 1, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 826.2604126278726, 7436.343713650853, 0.0324464383486897, 0.0324464383486897),
    @("FAPs", "Tgfb1", "Sdc2", "MuSCs", 3, 1, 18.65566266666667, 55.966988, 0.04786922362394307, 0.04786922362394307, 3, 1, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 356.4886202240783, 3208.397582016704, 0.01399895948218374, 0.01399895948218374),
    @("FAPs", "Tgfb1", "Sdc2", "Resolving-Mac", 3, 1, 18.65566266666667, 55.966988, 0.04786922362394307, 0.04786922362394307, 2, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 2.182445134168444, 19.642006207516, [double]"8.570248606002921E-05", [double]"8.570248606002922E-05"),
    @("Inflammatory-Mac", "Tgfb1", "Sdc2", "ECs", 3, 1, 171.0598806666667, 513.1796420000001, 0.4389285884413335, 0.4389285884413335, 3, 1, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 312.4521619897907, 2812.069457908116, 0.01226969083351502, 0.01226969083351502),
    @("Inflammatory-Mac", "Tgfb1", "Sdc2", "FAPs", 3, 1, 171.0598806666667, 513.1796420000001, 0.4389285884413335, 0.4389285884413335, 3, 1, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 7576.252321299548, 68186.27089169592, 0.2975120193346057, 0.2975120193346057),
    @("Inflammatory-Mac", "Tgfb1", "Sdc2", "MuSCs", 3, 1, 171.0598806666667, 513.1796420000001, 0.4389285884413335, 0.4389285884413335, 3, 1, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 3268.760907834927, 29418.84817051434, 0.12836104411121, 0.12836104411121),
    @("Inflammatory-Mac", "Tgfb1", "Sdc2", "Resolving-Mac", 3, 1, 171.0598806666667, 513.1796420000001, 0.4389285884413335, 0.4389285884413335, 2, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 20.01155417971044, 180.103987617394, 0.000785834162002711, 0.0007858341620027111),
    @("MuSCs", "Tgfb1", "Sdc2", "ECs", 3, 1, 12.628047, 37.884141, 0.0324027517316099, 0.0324027517316099, 3, 1, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 23.065961296602, 207.593651669418, 0.0009057777423744542, 0.0009057777423744543),
    @("MuSCs", "Tgfb1", "Sdc2", "FAPs", 3, 1, 12.628047, 37.884141, 0.0324027517316099, 0.0324027517316099, 3, 1, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 559.2969551034711, 5033.672595931239, 0.02196304445309023, 0.02196304445309023),
    @("MuSCs", "Tgfb1", "Sdc2", "MuSCs", 3, 1, 12.628047, 37.884141, 0.0324027517316099, 0.0324027517316099, 3, 1, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 241.307700058192, 2171.769300523728, 0.009475917390379056, 0.009475917390379056),
    @("MuSCs", "Tgfb1", "Sdc2", "Resolving-Mac", 3, 1, 12.628047, 37.884141, 0.0324027517316099, 0.0324027517316099, 2, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 1.477300496993, 13.295704472937, [double]"5.80121457661556E-05", [double]"5.801214576615561E-05"),
    @("Resolving-Mac", "Tgfb1", "Sdc2", "ECs", 3, 1, 139.5154473333333, 418.546342, 0.3579876130225438, 0.3579876130225438, 3, 1, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 254.8341725738573, 2293.507553164716, 0.01000708873762892, 0.01000708873762892),
    @("Resolving-Mac", "Tgfb1", "Sdc2", "FAPs", 3, 1, 139.5154473333333, 418.546342, 0.3579876130225438, 0.3579876130225438, 3, 1, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 6179.147486815024, 55612.32738133522, 0.242649078806467, 0.242649078806467),
    @("Resolving-Mac", "Tgfb1", "Sdc2", "MuSCs", 3, 1, 139.5154473333333, 418.546342, 0.3579876130225438, 0.3579876130225438, 3, 1, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 2665.982453074215, 23993.84207766793, 0.1046905236900406, 0.1046905236900406),
    @("Resolving-Mac", "Tgfb1", "Sdc2", "Resolving-Mac", 3, 1, 139.5154473333333, 418.546342, 0.3579876130225438, 0.3579876130225438, 2, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 16.32130761658822, 146.891768549294, 0.0006409217884073235, 0.0006409217884073236)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
